$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.433578968048096
$ws.Range("B1").Value = 3.567338943481445
$ws.Range("C1").Value = 5.327566146850586
$ws.Range("D1").Value = 1.730585813522339
$ws.Range("E1").Value = 0.9697878956794739
